$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "308.52"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "1.13%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "36.26"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "0.84%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.048"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "0.59%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.08133"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "0.85%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.992"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "5.47%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "4.149"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "0.00%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "7.861"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "0.33%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.9262"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-0.49%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1479"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "14.95%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1940"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "1.85%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.09088"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-1.34%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03524"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "0.34%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09869"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001418"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-0.07%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.006578"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "3.44%"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "5.10%"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "8.12%"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.3451"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "-0.01%"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.1313"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-2.42%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.810"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "-7.91%"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "-7.51%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04367"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-1.04%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001235"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "0.02%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004163"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-11.66%"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0001301"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "0.03%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02134"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "9.47%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.05122"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "-0.70%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007469"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-0.91%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.009987"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "-1.70%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1368"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-0.30%"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-1.81%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.009695"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-10.22%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006271"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-1.25%"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "0.02%"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "1.94%"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "-3.52%"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.02%"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "0.02%"
